$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range A1:N24 (with header row) ascending by column A (Year)
$rng = $ws.Range("A1:N24")
$rng.Sort($ws.Range("A1"), 1, $null, $null, 1, 0, 1, 1, $false, $null, $null, 1) | Out-Null

# Update the active selection to match the post-sort cursor position
$ws.Range("D10").Select() | Out-Null
